# NEMO related updates from shared google doc in pre ignored file.
# Adds 22 new rows (233, 235-254) of variable/comment data plus the
# associated shared strings, matching the upstream xlsx diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C233').Value = 'volcello'
$ws.Range('F233').Value = 'Not available in NEMO-OPA.'
$ws.Range('G233').Value = 'Raffaele Bernardello'
$ws.Range('I233').Value = 'grid-cell volume ca. 2000.'

$ws.Range('C235').Value = 'pso'
$ws.Range('F235').Value = 'Not available in NEMO.'

$ws.Range('C236').Value = 'msftmz'
$ws.Range('F236').Value = 'Not available in NEMO-OPA before CMIP6 starts. zomsflgo requires the subasins.nc file and a namelist parameter (Poleward Transport Diagnostic) / basin-wide variables are: zomsfatl zomsfpac zomsfind zomsfipc. NEMO-OPA - volume meridional stream function is available (zomsfglo). Maybe mass streamfunction can be obtained multiplying it by potential density (sea_water_sigma_theta) in the file_def xml ?'
$ws.Range('G236').Value = 'Etienne Tourigny'
$ws.Range('I236').Value = 'Overturning mass streamfunction arising from all advective mass transport processes, resolved and parameterized.'

$ws.Range('C237').Value = 'msftmrho'
$ws.Range('F237').Value = 'Not available in NEMO-OPA before CMIP6 starts. NEMO-OPA - volume meridional stream function is available (zomsfglo). Maybe mass streamfunction can be obtained multiplying it by potential density (sea_water_sigma_theta) in the file_def xml ?  '
$ws.Range('G237').Value = 'Raffaele Bernardello'
$ws.Range('I237').Value = 'Overturning mass streamfunction arising from all advective mass transport processes, resolved and parameterized.'

$ws.Range('C238').Value = 'msftyrho'
$ws.Range('F238').Value = 'Not available in NEMO-OPA before CMIP6 starts. NEMO-OPA - I guess it''s the same as above only rotated in case y does not align exactly with north-south direction.'
$ws.Range('G238').Value = 'Raffaele Bernardello'
$ws.Range('I238').Value = 'Overturning mass streamfunction arising from all advective mass transport processes, resolved and parameterized.'

$ws.Range('C239').Value = 'msftmzmpa'
$ws.Range('F239').Value = 'Not available in NEMO-OPA before CMIP6 starts.'
$ws.Range('G239').Value = 'Raffaele Bernardello'
$ws.Range('I239').Value = 'CMIP5 called this ''due to Bolus Advection''.  Name change respects the more general physics of the mesoscale parameterizations.'

$ws.Range('C240').Value = 'msftmrhompa'
$ws.Range('F240').Value = 'Not available in NEMO-OPA before CMIP6 starts.'
$ws.Range('G240').Value = 'Raffaele Bernardello'
$ws.Range('I240').Value = 'CMIP5 called this ''due to Bolus Advection''.  Name change respects the more general physics of the mesoscale parameterizations.'

$ws.Range('C241').Value = 'msftyzmpa'
$ws.Range('F241').Value = 'Not available in NEMO-OPA before CMIP6 starts.'
$ws.Range('G241').Value = 'Raffaele Bernardello'
$ws.Range('I241').Value = 'CMIP5 called this ''due to Bolus Advection''.  Name change respects the more general physics of the mesoscale parameterizations.'

$ws.Range('C242').Value = 'msftyrhompa'
$ws.Range('F242').Value = 'Not available in NEMO-OPA before CMIP6 starts.'
$ws.Range('G242').Value = 'Raffaele Bernardello'
$ws.Range('I242').Value = 'CMIP5 called this ''due to Bolus Advection''.  Name change respects the more general physics of the mesoscale parameterizations.'

$ws.Range('C243').Value = 'msftmzsmpa'
$ws.Range('F243').Value = 'Not available in NEMO-OPA before CMIP6 starts.'
$ws.Range('G243').Value = 'Raffaele Bernardello'
$ws.Range('I243').Value = 'Report only if there is a submesoscale eddy parameterization.'

$ws.Range('C244').Value = 'msftyzsmpa'
$ws.Range('F244').Value = 'Not available in NEMO-OPA before CMIP6 starts.'
$ws.Range('G244').Value = 'Raffaele Bernardello'
$ws.Range('I244').Value = 'Report only if there is a submesoscale eddy parameterization.'

$ws.Range('C245').Value = 'hfbasinpmdiff'
$ws.Range('F245').Value = 'Not available in NEMO-OPA before CMIP6 starts.'
$ws.Range('G245').Value = 'Raffaele Bernardello'
$ws.Range('I245').Value = 'Contributions to heat transport from parameterized mesoscale eddy-induced diffusive transport (i.e., neutral diffusion). Diagnosed here as a function of latitude and basin.'

$ws.Range('C246').Value = 'hfbasinpsmadv'
$ws.Range('F246').Value = 'Not available in NEMO-OPA before CMIP6 starts.'
$ws.Range('G246').Value = 'Raffaele Bernardello'
$ws.Range('I246').Value = 'Contributions to heat transport from parameterized mesoscale eddy-induced advective transport. Diagnosed here as a function of latitude and basin.  Use Celsius for temperature scale.'

$ws.Range('C247').Value = 'hfbasinpadv'
$ws.Range('F247').Value = 'Not available in NEMO-OPA before CMIP6 starts.'
$ws.Range('G247').Value = 'Raffaele Bernardello'
$ws.Range('I247').Value = 'Contributions to heat transport from parameterized eddy-induced advective transport due to any subgrid advective process. Diagnosed here as a function of latitude and basin.  Use Celsius for temperature scale.'

$ws.Range('C248').Value = 'wfcorr'
$ws.Range('F248').Value = 'Not available in NEMO-OPA.'
$ws.Range('G248').Value = 'Raffaele Bernardello'
$ws.Range('I248').Value = 'Positive flux implies correction adds water to ocean.'

$ws.Range('C249').Value = 'sfriver'
$ws.Range('F249').Value = 'Not available in NEMO-OPA, i.e. it makes no sence to make it availble because it is zero. It looks like it is assumed zero in NEMO, not 100% sure though.'
$ws.Range('G249').Value = 'Raffaele Bernardello'
$ws.Range('I249').Value = 'This field is physical, and it arises when rivers carry a nonzero salt content.  Often this is zero, with rivers assumed to be fresh.'

$ws.Range('C250').Value = 'hfsifrazil'
$ws.Range('F250').Value = 'Not available in NEMO-LIM, not in NEMO anywhere'
$ws.Range('G250').Value = 'Raffaele Bernardello'

$ws.Range('C251').Value = 'hfsifrazil2d'
$ws.Range('F251').Value = 'Not available in NEMO-LIM, not in NEMO anywhere'
$ws.Range('G251').Value = 'Raffaele Bernardello'

$ws.Range('C252').Value = 'hfcorr'
$ws.Range('F252').Value = 'Not available in NEMO-OPA.'
$ws.Range('G252').Value = 'Raffaele Bernardello'

$ws.Range('C253').Value = 'tauucorr'
$ws.Range('F253').Value = 'Not available in NEMO-OPA.'
$ws.Range('G253').Value = 'Raffaele Bernardello'
$ws.Range('I253').Value = 'This is the stress on the liquid ocean from overlying atmosphere, sea ice, ice shelf, etc.'

$ws.Range('C254').Value = 'tauvcorr'
$ws.Range('F254').Value = 'Not available in NEMO-OPA.'
$ws.Range('G254').Value = 'Raffaele Bernardello'
$ws.Range('I254').Value = 'This is the stress on the liquid ocean from overlying atmosphere, sea ice, ice shelf, etc.'

# Recreate the trailing "touched but empty" rows left behind at the
# bottom of the sheet (rows 1048567-1048576, row height 12.8) as seen
# in the authored workbook.
for ($r = 1048567; $r -le 1048576; $r++) {
    $ws.Rows.Item($r).RowHeight = 12.8
}

# Restore the cursor/selection to the cell the author left active (A234).
$ws.Range("A234").Select()

